# Update crypto price/volume figures to latest scraped values.
# Each cell is forced to Text format before assignment so that
# numeric-looking strings (e.g. '0.5243') are not silently
# converted to floating point numbers by Excel, then the cell
# style is reset to "Normal" so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '26.072.03'
Set-TextValue 'E2' '  -0.07%  '
Set-TextValue 'D3' '1.638.74'
Set-TextValue 'E3' '  -1.70%  '
Set-TextValue 'E4' '  -0.16%  '
Set-TextValue 'D5' '213.44'
Set-TextValue 'E5' '  +1.89%  '
Set-TextValue 'D6' '0.5243'
Set-TextValue 'E6' '  -0.06%  '
Set-TextValue 'D7' '1.002'
Set-TextValue 'E7' '  -0.13%  '
Set-TextValue 'E8' '  -0.87%  '
Set-TextValue 'D9' '0.06289'
Set-TextValue 'E9' '  +0.00%  '
Set-TextValue 'D10' '20.61'
Set-TextValue 'E10' '  -2.23%  '
Set-TextValue 'D11' '0.07678'
Set-TextValue 'E11' '  +1.94%  '
Set-TextValue 'D12' '1.625.61'
Set-TextValue 'E12' '  -2.53%  '
Set-TextValue 'D13' '4.407'
Set-TextValue 'E13' '  -0.50%  '
Set-TextValue 'D14' '1.862.23'
Set-TextValue 'E14' '  -1.72%  '
Set-TextValue 'D15' '0.5515'
Set-TextValue 'E15' '  +0.07%  '
Set-TextValue 'D16' '0.0₅8266'
Set-TextValue 'E16' '  +4.43%  '
Set-TextValue 'D17' '64.94'
Set-TextValue 'E17' '  -2.18%  '
Set-TextValue 'D18' '26.060.52'
Set-TextValue 'E18' '  -0.30%  '
Set-TextValue 'D19' '1.002'
Set-TextValue 'D20' '4.684'
Set-TextValue 'E20' '  -0.73%  '
Set-TextValue 'D21' '188.07'
Set-TextValue 'E21' '  +1.06%  '
Set-TextValue 'D22' '10.18'
Set-TextValue 'E22' '  -0.80%  '
Set-TextValue 'D23' '6.158'
Set-TextValue 'E23' '  +0.08%  '
Set-TextValue 'D24' '1.002'
Set-TextValue 'E24' '  -0.12%  '
Set-TextValue 'D25' '145.60'
Set-TextValue 'E25' '  -2.43%  '
Set-TextValue 'D26' '0.1214'
Set-TextValue 'E26' '  -2.18%  '
Set-TextValue 'D27' '7.410'
Set-TextValue 'E27' '  -0.69%  '
Set-TextValue 'D28' '15.83'
Set-TextValue 'E28' '  -0.31%  '
Set-TextValue 'D29' '1.399'
Set-TextValue 'E29' '  +3.61%  '
Set-TextValue 'D30' '0.05952'
Set-TextValue 'E30' '  -5.49%  '
Set-TextValue 'E31' '  -1.52%  '
Set-TextValue 'D32' '3.434'
Set-TextValue 'E32' '  -1.44%  '
Set-TextValue 'D33' '3.402'
Set-TextValue 'E33' '  -0.17%  '
Set-TextValue 'E34' '  +0.57%  '
Set-TextValue 'D35' '0.9824'
Set-TextValue 'E35' '  -1.54%  '
Set-TextValue 'D36' '2.396'
Set-TextValue 'E36' '  -0.42%  '
Set-TextValue 'D37' '2.758'
Set-TextValue 'E37' '  +1.10%  '
Set-TextValue 'D38' '0.5659'
Set-TextValue 'E38' '  -5.82%  '
Set-TextValue 'D39' '0.01616'
Set-TextValue 'E39' '  +0.17%  '
Set-TextValue 'D40' '0.8499'
Set-TextValue 'E40' '  -2.45%  '
Set-TextValue 'E41' '  -0.20%  '
Set-TextValue 'D42' '5.710'
Set-TextValue 'E42' '  -6.16%  '
Set-TextValue 'D43' '1.032.72'
Set-TextValue 'E43' '  -7.01%  '
Set-TextValue 'D44' '100.24'
Set-TextValue 'E44' '  +0.49%  '
Set-TextValue 'D45' '1.788.01'
Set-TextValue 'D46' '0.0₈110'
Set-TextValue 'E46' '  +0.64%  '
Set-TextValue 'D47' '55.79'
Set-TextValue 'E47' '  +0.92%  '
Set-TextValue 'D48' '1.002'
Set-TextValue 'E48' '  +0.07%  '
Set-TextValue 'D49' '8.078'
Set-TextValue 'E49' '  +0.37%  '
Set-TextValue 'D50' '0.05145'
Set-TextValue 'E50' '  -1.65%  '
Set-TextValue 'D51' '0.4216'
Set-TextValue 'E51' '  -0.67%  '
